$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 16, shifting existing row 16..115 down to 17..116.
$ws.Rows("16:16").Insert()

# Fill the new row 16 with the new weekly price-report entry.
$ws.Cells.Item(16, 1).Value = 11
$ws.Cells.Item(16, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(16, 3).Value = "Bíobío"
$ws.Cells.Item(16, 4).Value = 45035
$ws.Cells.Item(16, 5).Value = 8
$ws.Cells.Item(16, 6).Value = 100112012
$ws.Cells.Item(16, 7).Value = "Espinaca"
$ws.Cells.Item(16, 8).Value = "Sin especificar"
$ws.Cells.Item(16, 9).Value = "Primera"
$ws.Cells.Item(16, 10).Value = 40
$ws.Cells.Item(16, 11).Value = 10000
$ws.Cells.Item(16, 12).Value = 11000
$ws.Cells.Item(16, 13).Value = 10500
$ws.Cells.Item(16, 14).Value = "$/cuna 10 kilos"
$ws.Cells.Item(16, 15).Value = "Región Metropolitana"
$ws.Cells.Item(16, 16).Value = 1050
$ws.Cells.Item(16, 17).Value = 10
$ws.Cells.Item(16, 18).Value = "Hortaliza"
